$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = 0
$ws.Range("F24").Value = 1
$ws.Range("F26").Value = 1
$ws.Range("F28").Value = 0
$ws.Range("F30").Value = 2
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 1
$ws.Range("F33").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("F37").Value = 1
$ws.Range("F38").Value = 1
$ws.Range("F39").Value = 2
$ws.Range("F40").Value = 2
$ws.Range("F41").Value = 1
$ws.Range("F42").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("F47").Value = 2
$ws.Range("F49").Value = 0
$ws.Range("F52").Value = 0
$ws.Range("F55").Value = 0
$ws.Range("F57").Value = 0
$ws.Range("F58").Value = 1
$ws.Range("F59").Value = 0
$ws.Range("F60").Value = 1
$ws.Range("F61").Value = 0
$ws.Range("F62").Value = 1
$ws.Range("F63").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("F68").Value = 1
$ws.Range("F70").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("F72").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("F76").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("F78").Value = 1
$ws.Range("F79").Value = 2
$ws.Range("F80").Value = 0
$ws.Range("F81").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("F87").Value = 1
$ws.Range("F88").Value = 0
$ws.Range("F89").Value = 1
$ws.Range("F92").Value = 0
$ws.Range("F95").Value = 1
$ws.Range("F97").Value = 0
$ws.Range("F101").Value = 2
$ws.Range("F103").Value = 2
$ws.Range("F104").Value = 0
$ws.Range("F106").Value = 1
$ws.Range("F108").Value = 0
$ws.Range("F109").Value = 0
$ws.Range("F110").Value = 1
$ws.Range("F112").Value = 1
$ws.Range("F113").Value = 2
$ws.Range("F118").Value = 2
$ws.Range("F119").Value = 2
$ws.Range("F120").Value = 2
$ws.Range("F122").Value = 0
$ws.Range("F124").Value = 0
$ws.Range("F130").Value = 1
$ws.Range("F131").Value = 0
$ws.Range("F132").Value = 0
$ws.Range("F133").Value = 1
$ws.Range("F136").Value = 2
$ws.Range("F137").Value = 1
$ws.Range("F138").Value = 0
$ws.Range("F139").Value = 0
$ws.Range("F140").Value = 0
$ws.Range("F141").Value = 0
$ws.Range("F142").Value = 0
$ws.Range("F143").Value = 2
$ws.Range("F144").Value = 2
